$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model header columns
$ws.Range("E1").Value = "sdfcsfs"
$ws.Range("F1").Value = "adasd"

# Feature flags for new model columns
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1

$ws.Range("F6").Value = 1

$ws.Range("E8").Value = 1
